# Algeria Ligue 1 2023-2024 - script update (19-12-2023 02:45)
# - Corrects the home/away order for three "Round" fixtures whose scraped
#   rows had drifted out of sequence (rows 16-18 rotate, rows 35-36 swap).
# - Appends three newly scraped fixtures (Saoura-Magra, MC Alger-Khenchela,
#   Oran-Biskra) as rows 72-74.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-MatchRow {
    param($Row, $F,$G,$H,$I,$J,$K,$L,$M,$N,$O,$P,$Q,$R,$S,$T,$U,$V)

    $ws.Cells.Item($Row, 6).Value  = $F   # home
    $ws.Cells.Item($Row, 7).Value  = $G   # home_ft_gols
    $ws.Cells.Item($Row, 8).Value  = $H   # away
    $ws.Cells.Item($Row, 9).Value  = $I   # away_ft_gols
    $ws.Cells.Item($Row, 10).Value = $J   # home_opening_odds
    $ws.Cells.Item($Row, 11).Value = $K   # home_opening_data_hora
    $ws.Cells.Item($Row, 12).Value = $L   # home_closing_odds
    $ws.Cells.Item($Row, 13).Value = $M   # home_closing_data_hora
    $ws.Cells.Item($Row, 14).Value = $N   # draw_opening_odds
    $ws.Cells.Item($Row, 15).Value = $O   # draw_opening_data_hora
    $ws.Cells.Item($Row, 16).Value = $P   # draw_closing_odds
    $ws.Cells.Item($Row, 17).Value = $Q   # draw_closing_data_hora
    $ws.Cells.Item($Row, 18).Value = $R   # away_opening_odds
    $ws.Cells.Item($Row, 19).Value = $S   # away_opening_data_hora
    $ws.Cells.Item($Row, 20).Value = $T   # away_closing_odds
    $ws.Cells.Item($Row, 21).Value = $U   # away_closing_data_hora
    $ws.Cells.Item($Row, 22).Value = $V   # url_partida
}

# --- Rows 16-18: same matchday, rows had rotated out of order -------------
# Row 16 -> Khenchela vs Kabylie (was row 18's data)
Set-MatchRow 16 "Khenchela" 2 "Kabylie" 1 `
    2.63 "28/09/2023 04:12" 2.05 "29/09/2023 16:41" `
    2.62 "28/09/2023 04:12" 2.75 "29/09/2023 16:41" `
    3.02 "28/09/2023 04:12" 4.88 "29/09/2023 16:27" `
    "https://www.betexplorer.com/football/algeria/ligue-1/khenchela-kabylie/pUZYGLcr/"

# Row 17 -> Magra vs Biskra (was row 16's data)
Set-MatchRow 17 "Magra" 2 "Biskra" 1 `
    1.98 "28/09/2023 04:12" 2.1 "29/09/2023 16:44" `
    2.89 "28/09/2023 04:12" 2.77 "29/09/2023 16:44" `
    3.97 "28/09/2023 04:12" 4.54 "29/09/2023 16:44" `
    "https://www.betexplorer.com/football/algeria/ligue-1/magra-biskra/OKYxGuDl/"

# Row 18 -> US Souf vs Oran (was row 17's data)
Set-MatchRow 18 "US Souf" 0 "Oran" 0 `
    2.49 "28/09/2023 19:27" 2.14 "29/09/2023 13:29" `
    2.88 "28/09/2023 19:27" 2.74 "29/09/2023 14:49" `
    3.18 "28/09/2023 19:27" 4.43 "29/09/2023 15:47" `
    "https://www.betexplorer.com/football/algeria/ligue-1/us-souf-oran/6qOsFaSf/"

# --- Rows 35-36: swapped pair ---------------------------------------------
# Row 35 -> Khenchela vs Biskra (was row 36's data)
Set-MatchRow 35 "Khenchela" 0 "Biskra" 1 `
    1.74 "10/11/2023 03:13" 1.29 "11/11/2023 10:24" `
    3.19 "10/11/2023 03:13" 4.87 "11/11/2023 14:54" `
    4.9 "10/11/2023 03:13" 13.77 "11/11/2023 14:54" `
    "https://www.betexplorer.com/football/algeria/ligue-1/khenchela-biskra/GbL62yef/"

# Row 36 -> Magra vs Ben Aknoun (was row 35's data)
Set-MatchRow 36 "Magra" 3 "Ben Aknoun" 1 `
    1.61 "11/11/2023 10:12" 1.62 "11/11/2023 14:48" `
    3.51 "11/11/2023 10:12" 3.56 "11/11/2023 14:48" `
    5.97 "11/11/2023 10:12" 6.3 "11/11/2023 14:48" `
    "https://www.betexplorer.com/football/algeria/ligue-1/magra-es-ben-aknoun/lCJE0FP6/"

# --- New rows 72-74 (append, copying formatting from the last row) -------
$ws.Range("A71:V71").Copy()
$ws.Range("A72:V74").PasteSpecial(-4122)
$excel.CutCopyMode = $false

function Set-NewMatchRow {
    param($Row,$Idx,$Date,$F,$G,$H,$I,$J,$K,$L,$M,$N,$O,$P,$Q,$R,$S,$T,$U,$V)

    $ws.Cells.Item($Row, 1).Value  = $Idx
    $ws.Cells.Item($Row, 2).Value  = "algeria"
    $ws.Cells.Item($Row, 3).Value  = "ligue-1"
    $ws.Cells.Item($Row, 4).Value  = "2023-2024"
    $ws.Cells.Item($Row, 5).Value  = $Date

    Set-MatchRow $Row $F $G $H $I $J $K $L $M $N $O $P $Q $R $S $T $U $V
}

# Row 72 (index 71): Saoura 2 x 1 Magra
Set-NewMatchRow 72 71 45276.70833333334 "Saoura" 2 "Magra" 1 `
    1.42 "14/12/2023 09:42" 1.39 "16/12/2023 16:02" `
    3.85 "14/12/2023 09:42" 4.18 "16/12/2023 15:02" `
    8.01 "14/12/2023 09:42" 10.57 "16/12/2023 16:02" `
    "https://www.betexplorer.com/football/algeria/ligue-1/saoura-magra/tWNi7Z35/"

# Row 73 (index 72): MC Alger 3 x 0 Khenchela
Set-NewMatchRow 73 72 45276.75 "MC Alger" 3 "Khenchela" 0 `
    1.42 "14/12/2023 09:42" 1.33 "16/12/2023 17:24" `
    3.86 "14/12/2023 09:42" 4.52 "16/12/2023 17:24" `
    7.89 "14/12/2023 09:42" 12.14 "16/12/2023 17:24" `
    "https://www.betexplorer.com/football/algeria/ligue-1/mc-alger-khenchela/p4Ur9eZh/"

# Row 74 (index 73): Oran 0 x 1 Biskra
Set-NewMatchRow 74 73 45276.79166666666 "Oran" 0 "Biskra" 1 `
    2.12 "14/12/2023 09:42" 1.83 "16/12/2023 18:55" `
    2.79 "14/12/2023 09:42" 3.16 "16/12/2023 18:59" `
    3.8 "14/12/2023 09:42" 5.16 "16/12/2023 18:59" `
    "https://www.betexplorer.com/football/algeria/ligue-1/oran-biskra/6DVvAyKn/"
